# "Generate Report for Handoff"
#
# ebae9fcd-9396-40e3-b465-4c773446c4db was handed off again, making it the
# most recently handed-off file. The localization-status report re-orders
# it ahead of 7c4eecda-dce9-469c-acc9-eb8e85e95c20 (their rows effectively
# swap) on every sheet, and the "Latest Handoff Datetime" column (shared
# across the whole per-language table) is refreshed to the new handoff
# timestamp.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Overview": A2/A3 (file name + hyperlink) swap; B/C unchanged.
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value = "ebae9fcd-9396-40e3-b465-4c773446c4db.md"
$wsOverview.Range("A3").Value = "7c4eecda-dce9-469c-acc9-eb8e85e95c20.md"

$wsOverview.Hyperlinks.Delete()

$overviewLinks = @(
  @{Cell="A2"; Url="https://github.com/OpenLocalizationTest/oltest/blob/7600a27c7b56ab9ade92818153d0bc2ebb55e51e/e2e/7c4eecda-dce9-469c-acc9-eb8e85e95c20.md"; Text="7c4eecda-dce9-469c-acc9-eb8e85e95c20.md"},
  @{Cell="A3"; Url="https://github.com/OpenLocalizationTest/oltest/blob/7600a27c7b56ab9ade92818153d0bc2ebb55e51e/e2e/ebae9fcd-9396-40e3-b465-4c773446c4db.md"; Text="ebae9fcd-9396-40e3-b465-4c773446c4db.md"},
  @{Cell="A4"; Url="https://github.com/OpenLocalizationTest/oltest/blob/7600a27c7b56ab9ade92818153d0bc2ebb55e51e/e2e/41fa2a1e-e5f6-419a-8cae-3684c8394aac.md"; Text="41fa2a1e-e5f6-419a-8cae-3684c8394aac.md"},
  @{Cell="A5"; Url="https://github.com/OpenLocalizationTest/oltest/blob/7600a27c7b56ab9ade92818153d0bc2ebb55e51e/e2e/f491a28a-ae0e-4d0d-98aa-0ad501f29e48.md"; Text="f491a28a-ae0e-4d0d-98aa-0ad501f29e48.md"},
  @{Cell="A6"; Url="https://github.com/OpenLocalizationTest/oltest/blob/7600a27c7b56ab9ade92818153d0bc2ebb55e51e/.localization-config"; Text=".localization-config"}
)
foreach ($l in $overviewLinks) {
    $wsOverview.Hyperlinks.Add($wsOverview.Range($l.Cell), $l.Url, "", "", $l.Text) | Out-Null
}

# ---------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("A2").Value = "ebae9fcd-9396-40e3-b465-4c773446c4db.md"
$wsZh.Range("C2").Value = "ebae9fcd-9396-40e3-b465-4c773446c4db.7ace181aef5501e3eca4e1cef8bdf56f0a6caf01.zh-cn.xlf"
$wsZh.Range("D2").Value = "2016-03-03 15:56:58"
$wsZh.Range("E2").Value = "ebae9fcd-9396-40e3-b465-4c773446c4db.md"
$wsZh.Range("F2").Value = "ebae9fcd-9396-40e3-b465-4c773446c4db.7ace181aef5501e3eca4e1cef8bdf56f0a6caf01.zh-cn.xlf"

$wsZh.Range("A3").Value = "7c4eecda-dce9-469c-acc9-eb8e85e95c20.md"
$wsZh.Range("C3").Value = "7c4eecda-dce9-469c-acc9-eb8e85e95c20.0bd6b73cc9e1abe76e3214de4ab6371dbbe65a8d.zh-cn.xlf"
$wsZh.Range("D3").Value = "2016-03-03 15:56:58"
$wsZh.Range("E3").Value = "7c4eecda-dce9-469c-acc9-eb8e85e95c20.md"
$wsZh.Range("F3").Value = "7c4eecda-dce9-469c-acc9-eb8e85e95c20.0bd6b73cc9e1abe76e3214de4ab6371dbbe65a8d.zh-cn.xlf"

# Latest Handoff Datetime column is the same value across the whole sheet.
$wsZh.Range("D4").Value = "2016-03-03 15:56:58"
$wsZh.Range("D5").Value = "2016-03-03 15:56:58"

$wsZh.Hyperlinks.Delete()

$zhLinks = @(
  @{Cell="A2"; Url="https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/979cd360471c00fc3666f35959cb21421f2d4a06/e2e/ebae9fcd-9396-40e3-b465-4c773446c4db.md"; Text="ebae9fcd-9396-40e3-b465-4c773446c4db.md"},
  @{Cell="C2"; Url="https://github.com/OpenLocalizationTestOrg/olhandoff/blob/486ae0044df4aa399991f1eca667a2cd8bc45615/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/ebae9fcd-9396-40e3-b465-4c773446c4db.7ace181aef5501e3eca4e1cef8bdf56f0a6caf01.zh-cn.xlf"; Text="ebae9fcd-9396-40e3-b465-4c773446c4db.7ace181aef5501e3eca4e1cef8bdf56f0a6caf01.zh-cn.xlf"},
  @{Cell="E2"; Url="https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/979cd360471c00fc3666f35959cb21421f2d4a06/e2e/ebae9fcd-9396-40e3-b465-4c773446c4db.md"; Text="ebae9fcd-9396-40e3-b465-4c773446c4db.md"},
  @{Cell="F2"; Url="https://github.com/OpenLocalizationTestOrg/olhandback/blob/d609068120e679d10eba1360dcd6f85da1bddbc1/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/mt/ebae9fcd-9396-40e3-b465-4c773446c4db.7ace181aef5501e3eca4e1cef8bdf56f0a6caf01.zh-cn.xlf"; Text="ebae9fcd-9396-40e3-b465-4c773446c4db.7ace181aef5501e3eca4e1cef8bdf56f0a6caf01.zh-cn.xlf"},

  @{Cell="A3"; Url="https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/979cd360471c00fc3666f35959cb21421f2d4a06/e2e/7c4eecda-dce9-469c-acc9-eb8e85e95c20.md"; Text="7c4eecda-dce9-469c-acc9-eb8e85e95c20.md"},
  @{Cell="C3"; Url="https://github.com/OpenLocalizationTestOrg/olhandoff/blob/486ae0044df4aa399991f1eca667a2cd8bc45615/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/7c4eecda-dce9-469c-acc9-eb8e85e95c20.0bd6b73cc9e1abe76e3214de4ab6371dbbe65a8d.zh-cn.xlf"; Text="7c4eecda-dce9-469c-acc9-eb8e85e95c20.0bd6b73cc9e1abe76e3214de4ab6371dbbe65a8d.zh-cn.xlf"},
  @{Cell="E3"; Url="https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/979cd360471c00fc3666f35959cb21421f2d4a06/e2e/7c4eecda-dce9-469c-acc9-eb8e85e95c20.md"; Text="7c4eecda-dce9-469c-acc9-eb8e85e95c20.md"},
  @{Cell="F3"; Url="https://github.com/OpenLocalizationTestOrg/olhandback/blob/d609068120e679d10eba1360dcd6f85da1bddbc1/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/mt/7c4eecda-dce9-469c-acc9-eb8e85e95c20.0bd6b73cc9e1abe76e3214de4ab6371dbbe65a8d.zh-cn.xlf"; Text="7c4eecda-dce9-469c-acc9-eb8e85e95c20.0bd6b73cc9e1abe76e3214de4ab6371dbbe65a8d.zh-cn.xlf"},

  @{Cell="A4"; Url="https://github.com/OpenLocalizationTest/oltest/blob/7600a27c7b56ab9ade92818153d0bc2ebb55e51e/e2e/41fa2a1e-e5f6-419a-8cae-3684c8394aac.md"; Text="41fa2a1e-e5f6-419a-8cae-3684c8394aac.md"},
  @{Cell="C4"; Url="https://github.com/OpenLocalizationTestOrg/olhandoff/blob/486ae0044df4aa399991f1eca667a2cd8bc45615/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/41fa2a1e-e5f6-419a-8cae-3684c8394aac.e5e53929549bef51e53423cb6de6c93f89a0399f.zh-cn.xlf"; Text="41fa2a1e-e5f6-419a-8cae-3684c8394aac.e5e53929549bef51e53423cb6de6c93f89a0399f.zh-cn.xlf"},

  @{Cell="A5"; Url="https://github.com/OpenLocalizationTest/oltest/blob/7600a27c7b56ab9ade92818153d0bc2ebb55e51e/e2e/f491a28a-ae0e-4d0d-98aa-0ad501f29e48.md"; Text="f491a28a-ae0e-4d0d-98aa-0ad501f29e48.md"},
  @{Cell="C5"; Url="https://github.com/OpenLocalizationTestOrg/olhandoff/blob/486ae0044df4aa399991f1eca667a2cd8bc45615/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/f491a28a-ae0e-4d0d-98aa-0ad501f29e48.1755ebbbfef550e4347980f9ae77e572d2349b51.zh-cn.xlf"; Text="f491a28a-ae0e-4d0d-98aa-0ad501f29e48.1755ebbbfef550e4347980f9ae77e572d2349b51.zh-cn.xlf"},

  @{Cell="A6"; Url="https://github.com/OpenLocalizationTest/oltest/blob/7600a27c7b56ab9ade92818153d0bc2ebb55e51e/.localization-config"; Text=".localization-config"}
)
foreach ($l in $zhLinks) {
    $wsZh.Hyperlinks.Add($wsZh.Range($l.Cell), $l.Url, "", "", $l.Text) | Out-Null
}

# ---------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("A2").Value = "ebae9fcd-9396-40e3-b465-4c773446c4db.md"
$wsDe.Range("C2").Value = "ebae9fcd-9396-40e3-b465-4c773446c4db.7ace181aef5501e3eca4e1cef8bdf56f0a6caf01.de-de.xlf"
$wsDe.Range("D2").Value = "2016-03-03 15:57:12"
$wsDe.Range("E2").Value = "ebae9fcd-9396-40e3-b465-4c773446c4db.md"
$wsDe.Range("F2").Value = "ebae9fcd-9396-40e3-b465-4c773446c4db.7ace181aef5501e3eca4e1cef8bdf56f0a6caf01.de-de.xlf"

$wsDe.Range("A3").Value = "7c4eecda-dce9-469c-acc9-eb8e85e95c20.md"
$wsDe.Range("C3").Value = "7c4eecda-dce9-469c-acc9-eb8e85e95c20.0bd6b73cc9e1abe76e3214de4ab6371dbbe65a8d.de-de.xlf"
$wsDe.Range("D3").Value = "2016-03-03 15:57:12"
$wsDe.Range("E3").Value = "7c4eecda-dce9-469c-acc9-eb8e85e95c20.md"
$wsDe.Range("F3").Value = "7c4eecda-dce9-469c-acc9-eb8e85e95c20.0bd6b73cc9e1abe76e3214de4ab6371dbbe65a8d.de-de.xlf"

# Latest Handoff Datetime column is the same value across the whole sheet.
$wsDe.Range("D4").Value = "2016-03-03 15:57:12"
$wsDe.Range("D5").Value = "2016-03-03 15:57:12"

$wsDe.Hyperlinks.Delete()

$deLinks = @(
  @{Cell="A2"; Url="https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/bed7ad7e31e7313cd448918b4b20e19e8f0e4b4e/e2e/ebae9fcd-9396-40e3-b465-4c773446c4db.md"; Text="ebae9fcd-9396-40e3-b465-4c773446c4db.md"},
  @{Cell="C2"; Url="https://github.com/OpenLocalizationTestOrg/olhandoff/blob/57579b8ca32807529ef93e2a5474b5608833b03d/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/ebae9fcd-9396-40e3-b465-4c773446c4db.7ace181aef5501e3eca4e1cef8bdf56f0a6caf01.de-de.xlf"; Text="ebae9fcd-9396-40e3-b465-4c773446c4db.7ace181aef5501e3eca4e1cef8bdf56f0a6caf01.de-de.xlf"},
  @{Cell="E2"; Url="https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/bed7ad7e31e7313cd448918b4b20e19e8f0e4b4e/e2e/ebae9fcd-9396-40e3-b465-4c773446c4db.md"; Text="ebae9fcd-9396-40e3-b465-4c773446c4db.md"},
  @{Cell="F2"; Url="https://github.com/OpenLocalizationTestOrg/olhandback/blob/5bc2af341d03d2cb1b4d71a04ff415d0dc87b215/ol-handback/OpenLocalizationTestOrg/oltest.de-de/xinjiang/mt/ebae9fcd-9396-40e3-b465-4c773446c4db.7ace181aef5501e3eca4e1cef8bdf56f0a6caf01.de-de.xlf"; Text="ebae9fcd-9396-40e3-b465-4c773446c4db.7ace181aef5501e3eca4e1cef8bdf56f0a6caf01.de-de.xlf"},

  @{Cell="A3"; Url="https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/bed7ad7e31e7313cd448918b4b20e19e8f0e4b4e/e2e/7c4eecda-dce9-469c-acc9-eb8e85e95c20.md"; Text="7c4eecda-dce9-469c-acc9-eb8e85e95c20.md"},
  @{Cell="C3"; Url="https://github.com/OpenLocalizationTestOrg/olhandoff/blob/57579b8ca32807529ef93e2a5474b5608833b03d/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/7c4eecda-dce9-469c-acc9-eb8e85e95c20.0bd6b73cc9e1abe76e3214de4ab6371dbbe65a8d.de-de.xlf"; Text="7c4eecda-dce9-469c-acc9-eb8e85e95c20.0bd6b73cc9e1abe76e3214de4ab6371dbbe65a8d.de-de.xlf"},
  @{Cell="E3"; Url="https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/bed7ad7e31e7313cd448918b4b20e19e8f0e4b4e/e2e/7c4eecda-dce9-469c-acc9-eb8e85e95c20.md"; Text="7c4eecda-dce9-469c-acc9-eb8e85e95c20.md"},
  @{Cell="F3"; Url="https://github.com/OpenLocalizationTestOrg/olhandback/blob/5bc2af341d03d2cb1b4d71a04ff415d0dc87b215/ol-handback/OpenLocalizationTestOrg/oltest.de-de/xinjiang/mt/7c4eecda-dce9-469c-acc9-eb8e85e95c20.0bd6b73cc9e1abe76e3214de4ab6371dbbe65a8d.de-de.xlf"; Text="7c4eecda-dce9-469c-acc9-eb8e85e95c20.0bd6b73cc9e1abe76e3214de4ab6371dbbe65a8d.de-de.xlf"},

  @{Cell="A4"; Url="https://github.com/OpenLocalizationTest/oltest/blob/7600a27c7b56ab9ade92818153d0bc2ebb55e51e/e2e/41fa2a1e-e5f6-419a-8cae-3684c8394aac.md"; Text="41fa2a1e-e5f6-419a-8cae-3684c8394aac.md"},
  @{Cell="C4"; Url="https://github.com/OpenLocalizationTestOrg/olhandoff/blob/57579b8ca32807529ef93e2a5474b5608833b03d/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/41fa2a1e-e5f6-419a-8cae-3684c8394aac.e5e53929549bef51e53423cb6de6c93f89a0399f.de-de.xlf"; Text="41fa2a1e-e5f6-419a-8cae-3684c8394aac.e5e53929549bef51e53423cb6de6c93f89a0399f.de-de.xlf"},

  @{Cell="A5"; Url="https://github.com/OpenLocalizationTest/oltest/blob/7600a27c7b56ab9ade92818153d0bc2ebb55e51e/e2e/f491a28a-ae0e-4d0d-98aa-0ad501f29e48.md"; Text="f491a28a-ae0e-4d0d-98aa-0ad501f29e48.md"},
  @{Cell="C5"; Url="https://github.com/OpenLocalizationTestOrg/olhandoff/blob/57579b8ca32807529ef93e2a5474b5608833b03d/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/f491a28a-ae0e-4d0d-98aa-0ad501f29e48.1755ebbbfef550e4347980f9ae77e572d2349b51.de-de.xlf"; Text="f491a28a-ae0e-4d0d-98aa-0ad501f29e48.1755ebbbfef550e4347980f9ae77e572d2349b51.de-de.xlf"},

  @{Cell="A6"; Url="https://github.com/OpenLocalizationTest/oltest/blob/7600a27c7b56ab9ade92818153d0bc2ebb55e51e/.localization-config"; Text=".localization-config"}
)
foreach ($l in $deLinks) {
    $wsDe.Hyperlinks.Add($wsDe.Range($l.Cell), $l.Url, "", "", $l.Text) | Out-Null
}

Write-Host "Localization status report regenerated for handoff."
